# Study Plan Optimizer: add "Earliest Start Week" column to Activities sheet
# and refresh selections / column width per author's commit.

$wb = $excel.ActiveWorkbook

$activities = $wb.Worksheets.Item("Activities")
$weekAvail  = $wb.Worksheets.Item("Week Availability")

# --- Activities sheet: new "Earliest Start Week" column (E) ---

# Header cell: copy the formatting of the neighbouring header (D1) so the
# new header gets the same bold/centered style, then set its text.
$activities.Range("D1").Copy() | Out-Null
$activities.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$activities.Range("E1").Value = "Earliest Start Week"

# Data values for the new column.
$earliestStartWeeks = @(1, 2, 2, 8, 4, 1)
for ($i = 0; $i -lt $earliestStartWeeks.Length; $i++) {
    $row = 2 + $i
    $activities.Cells.Item($row, 5).Value = $earliestStartWeeks[$i]
}

$activities.Range("A1").Application.CutCopyMode = $false

# --- Week Availability sheet: widen column D slightly, refresh selection ---

$weekAvail.Columns("D").ColumnWidth = 11.17
$weekAvail.Range("J14").Select() | Out-Null

# --- Restore Activities as the active sheet/selection ---

$activities.Activate()
$activities.Range("A2").Select() | Out-Null
